# Update crypto price ticker data (auto-refresh snapshot).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '29.101.38'
$ws.Range("E2").Value = '  -1.68%  '

# Row 3
$ws.Range("D3").Value = '1.836.71'
$ws.Range("E3").Value = '  -1.20%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9990'
$ws.Range("E4").Value = '  -0.06%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.53'
$ws.Range("E5").Value = '  -1.74%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6765'
$ws.Range("E6").Value = '  -2.59%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9999'
$ws.Range("E7").Value = '  -0.02%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2980'
$ws.Range("E8").Value = '  -2.58%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07457'
$ws.Range("E9").Value = '  -3.09%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.06'
$ws.Range("E10").Value = '  -2.76%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07667'
$ws.Range("E11").Value = '  -1.19%  '

# Row 12
$ws.Range("D12").Value = '1.833.89'
$ws.Range("E12").Value = '  -1.15%  '

# Row 13
$ws.Range("E13").Value = '  -2.59%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6763'
$ws.Range("E14").Value = '  -2.18%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '86.43'
$ws.Range("E15").Value = '  -5.85%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.162'
$ws.Range("E16").Value = '  -6.20%  '

# Row 17
$ws.Range("D17").Value = '29.089.92'
$ws.Range("E17").Value = '  -1.72%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008260'
$ws.Range("E18").Value = '  -0.36%  '

# Row 19
$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '227.87'
$ws.Range("E19").Value = '  -5.18%  '

# Row 20
$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.49'
$ws.Range("E20").Value = '  -2.18%  '

# Row 21
$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9992'
$ws.Range("E21").Value = '  -0.10%  '

# Row 22
$ws.Range("B22").Value = 'Chainlink'
$ws.Range("C22").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.343'
$ws.Range("E22").Value = '  -3.34%  '

# Row 23
$ws.Range("B23").Value = 'BinanceUSD'
$ws.Range("C23").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.000'
$ws.Range("E23").Value = '  -0.03%  '

# Row 24
$ws.Range("B24").Value = 'Monero'
$ws.Range("C24").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '160.76'
$ws.Range("E24").Value = '  +0.42%  '

# Row 25
$ws.Range("B25").Value = 'Stellar'
$ws.Range("C25").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1432'
$ws.Range("E25").Value = '  -4.40%  '

# Row 26
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.696'
$ws.Range("E26").Value = '  -2.54%  '

# Row 27
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.01'
$ws.Range("E27").Value = '  -1.50%  '

# Row 28
$ws.Range("B28").Value = 'PancakeSwap'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.502'
$ws.Range("E28").Value = '  -2.15%  '

# Row 29
$ws.Range("B29").Value = 'Filecoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.247'
$ws.Range("E29").Value = '  -0.10%  '

# Row 30
$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.130'
$ws.Range("E30").Value = '  -1.32%  '

# Row 31
$ws.Range("B31").Value = 'Hedera'
$ws.Range("C31").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.05427'
$ws.Range("E31").Value = '  +6.56%  '

# Row 32
$ws.Range("B32").Value = 'Toncoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.197'
$ws.Range("E32").Value = '  -0.34%  '

# Row 33
$ws.Range("B33").Value = 'LidoDAOToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.862'
$ws.Range("E33").Value = '  -1.65%  '

# Row 34
$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7497'
$ws.Range("E34").Value = '  -2.93%  '

# Row 35
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.128'
$ws.Range("E35").Value = '  -2.04%  '

# Row 36
$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.682'
$ws.Range("E36").Value = '  -0.01%  '

# Row 37
$ws.Range("B37").Value = 'Maker'
$ws.Range("C37").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D37").Value = '1.303.13'
$ws.Range("E37").Value = '  -2.36%  '

# Row 38
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01809'
$ws.Range("E38").Value = '  -3.39%  '

# Row 39
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.709'
$ws.Range("E39").Value = '  -0.63%  '

# Row 40
$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9353'
$ws.Range("E40").Value = '  -3.02%  '

# Row 41
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.094'
$ws.Range("E41").Value = '  +5.22%  '

# Row 42
$ws.Range("B42").Value = 'Quant'
$ws.Range("C42").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '104.55'
$ws.Range("E42").Value = '  -2.06%  '

# Row 43
$ws.Range("B43").Value = 'PaxDollar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9990'
$ws.Range("E43").Value = '  -0.08%  '

# Row 44
$ws.Range("B44").Value = 'XinFinNetwork'
$ws.Range("C44").Value = 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.08232'
$ws.Range("E44").Value = '  +29.08%  '

# Row 45
$ws.Range("B45").Value = 'RocketPoolETH'
$ws.Range("C45").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D45").Value = '1.976.46'
$ws.Range("E45").Value = '  -1.30%  '

# Row 46
$ws.Range("B46").Value = 'Mantle'
$ws.Range("C46").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5178'
$ws.Range("E46").Value = '  -0.71%  '

# Row 47
$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000122'
$ws.Range("E47").Value = '  -3.42%  '

# Row 48
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.480'
$ws.Range("E48").Value = '  -3.02%  '

# Row 49
$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.761'
$ws.Range("E49").Value = '  -0.67%  '

# Row 50
$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '63.66'
$ws.Range("E50").Value = '  +0.17%  '

# Row 51
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05930'
$ws.Range("E51").Value = '  +0.10%  '

